# RMPP_PCOM7E_Skill_Matrix.xlsx
# This module was previously used for "SEPM" (Software Engineering Project
# Management) and is being re-used/re-uploaded for "RMPP". Update the
# "Evidence" column (E) text so it refers to RMPP instead of SEPM, tidy a
# couple of entries, and move the selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "RMPP - End of module assignment, Work experience - working  with many teams, third parties and  non-technical departments. Worked on design documentation."
$ws.Range("E3").Value = "RMPP - End of module assignment. Work experience. Working with different document types. Examples include combining diagrams and plans in the Systems Design Assignment, along with writing a manual for the implementation assignment."
$ws.Range("E4").Value = "participation in discussion forums"
$ws.Range("E5").Value = "participation in discussion forums"
$ws.Range("E9").Value = "RMPP-End of module assignment, discussion relating to research methods."
$ws.Range("E10").Value = "RMPP-End of module assignment, discussion relating to research methods.."
$ws.Range("E11").Value = "RMPP-End of module assignment, discussion relating to research methods."
$ws.Range("E12").Value = "RMPP-End of module assignment."
$ws.Range("E13").Value = "RMPP-End of module assignment."

# Move the saved view: scroll down a bit and land the selection on E13.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
